$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 413.5
$ws.Range("I6").Value = 85.166664
$ws.Range("J6").Value = 1398.5
$ws.Range("K6").Value = 255.499992
$ws.Range("L6").Value = 4195.5
$ws.Range("M6").Value = -143.499992
$ws.Range("N6").Value = -4419.5
$ws.Range("H17").Value = 954.0175
$ws.Range("J17").Value = 954.0175
$ws.Range("L17").Value = 2862.0525
$ws.Range("N17").Value = -3198.0525
$ws.Range("H28").Value = 2379740
$ws.Range("J28").Value = 12333
$ws.Range("L28").Value = 12333
$ws.Range("N28").Value = -13303
$ws.Range("H45").Value = 1900
$ws.Range("J45").Value = 0
$ws.Range("L45").Value = 0
$ws.Range("N45").ClearContents()
$ws.Range("H76").Value = 2607505.5
$ws.Range("I76").Value = 7810008
$ws.Range("J76").Value = 6254.1665
$ws.Range("K76").Value = 7810008
$ws.Range("L76").Value = 6254.1665
$ws.Range("M76").Value = -7809693
$ws.Range("N76").Value = -6884.1665
$ws.Range("H79").Value = 2607505.5
$ws.Range("I79").Value = 7810008
$ws.Range("J79").Value = 6254.1665
$ws.Range("K79").Value = 7810008
$ws.Range("L79").Value = 6254.1665
$ws.Range("M79").Value = -7808916
$ws.Range("N79").Value = -8438.166499999999
$ws.Range("H116").Value = 17225
$ws.Range("I116").Value = 51900
$ws.Range("J116").Value = 5666.6665
$ws.Range("K116").Value = 51900
$ws.Range("L116").Value = 5666.6665
$ws.Range("M116").Value = -48458
$ws.Range("N116").Value = -12550.6665
$ws.Range("H137").Value = 1259.6765
$ws.Range("J137").Value = 1600
$ws.Range("L137").Value = 4800
$ws.Range("N137").Value = -9900
$ws.Range("H138").Value = 2892.151
$ws.Range("I138").Value = 2623.926
$ws.Range("J138").Value = 3170.6924
$ws.Range("K138").Value = 7871.778
$ws.Range("L138").Value = 9512.0772
$ws.Range("M138").Value = -2731.778
$ws.Range("N138").Value = -19792.0772
$ws.Range("H141").Value = 1079211.1
$ws.Range("I141").Value = 1334833.8
$ws.Range("K141").Value = 4004501.4
$ws.Range("M141").Value = -3999321.4

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 1899
$ws.Range("I45").Value = 0
$ws.Range("K45").Value = 0
$ws.Range("M45").ClearContents()
$ws.Range("H105").Value = 0
$ws.Range("J105").Value = 0
$ws.Range("L105").Value = 0
$ws.Range("N105").ClearContents()
$ws.Range("H132").Value = 2286.5
$ws.Range("I132").Value = 1924.1333
$ws.Range("J132").Value = 2780.6365
$ws.Range("K132").Value = 5772.3999
$ws.Range("L132").Value = 8341.9095
$ws.Range("M132").Value = -3242.3999
$ws.Range("N132").Value = -13401.9095

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 2784.5715
$ws.Range("I107").Value = 2784.5715
$ws.Range("K107").Value = 2784.5715
$ws.Range("M107").Value = -864.5715

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1809.9259
$ws.Range("I31").Value = 1057.5883
$ws.Range("J31").Value = 3088.9
$ws.Range("K31").Value = 1057.5883
$ws.Range("L31").Value = 3088.9
$ws.Range("M31").Value = -762.5882999999999
$ws.Range("N31").Value = -3678.9
$ws.Range("H34").Value = 1809.9259
$ws.Range("I34").Value = 1057.5883
$ws.Range("J34").Value = 3088.9
$ws.Range("K34").Value = 1057.5883
$ws.Range("L34").Value = 3088.9
$ws.Range("M34").Value = -855.5882999999999
$ws.Range("N34").Value = -3492.9
$ws.Range("H134").Value = 1054.3334
$ws.Range("I134").Value = 1054.3334
$ws.Range("J134").Value = 0
$ws.Range("K134").Value = 3163.0002
$ws.Range("L134").Value = 0
$ws.Range("M134").Value = -628.0001999999999
$ws.Range("N134").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H7").Value = 2369.7
$ws.Range("I7").Value = 3416.5
$ws.Range("J7").Value = 799.5
$ws.Range("K7").Value = 10249.5
$ws.Range("L7").Value = 2398.5
$ws.Range("M7").Value = -10137.5
$ws.Range("N7").Value = -2622.5
$ws.Range("H38").Value = 431
$ws.Range("I38").Value = 58.75
$ws.Range("J38").Value = 927.3333
$ws.Range("K38").Value = 176.25
$ws.Range("L38").Value = 2781.9999
$ws.Range("M38").Value = 170.75
$ws.Range("N38").Value = -3475.9999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 3393.25
$ws.Range("I80").Value = 3393.25
$ws.Range("K80").Value = 3393.25
$ws.Range("M80").Value = -2395.25
$ws.Range("H83").Value = 3393.25
$ws.Range("I83").Value = 3393.25
$ws.Range("K83").Value = 16966.25
$ws.Range("M83").Value = -11974.25
$ws.Range("H132").Value = 2961031.8
$ws.Range("I132").Value = 4809329.5
$ws.Range("K132").Value = 14427988.5
$ws.Range("M132").Value = -14425458.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 5902
$ws.Range("I7").Value = 5902
$ws.Range("J7").Value = 0
$ws.Range("K7").Value = 5902
$ws.Range("L7").Value = 0
$ws.Range("M7").Value = -5790
$ws.Range("N7").ClearContents()
$ws.Range("H36").Value = 0
$ws.Range("J36").Value = 0
$ws.Range("L36").Value = 0
$ws.Range("N36").ClearContents()
$ws.Range("H40").Value = 5349.091
$ws.Range("I40").Value = 1730
$ws.Range("K40").Value = 1730
$ws.Range("M40").Value = -1594
$ws.Range("H68").Value = 2616.5454
$ws.Range("I68").Value = 2498.7778
$ws.Range("J68").Value = 3146.5
$ws.Range("K68").Value = 2498.7778
$ws.Range("L68").Value = 3146.5
$ws.Range("M68").Value = -1749.7778
$ws.Range("N68").Value = -4644.5
$ws.Range("H71").Value = 2616.5454
$ws.Range("I71").Value = 2498.7778
$ws.Range("J71").Value = 3146.5
$ws.Range("K71").Value = 12493.889
$ws.Range("L71").Value = 15732.5
$ws.Range("M71").Value = -8749.888999999999
$ws.Range("N71").Value = -23220.5
$ws.Range("H81").Value = 0
$ws.Range("J81").Value = 0
$ws.Range("L81").Value = 0
$ws.Range("N81").ClearContents()
$ws.Range("H82").Value = 2455
$ws.Range("I82").Value = 1588.6666
$ws.Range("J82").Value = 3494.6
$ws.Range("K82").Value = 1588.6666
$ws.Range("L82").Value = 3494.6
$ws.Range("M82").Value = -1227.6666
$ws.Range("N82").Value = -4216.6
$ws.Range("H84").Value = 0
$ws.Range("J84").Value = 0
$ws.Range("L84").Value = 0
$ws.Range("N84").ClearContents()
$ws.Range("H85").Value = 2455
$ws.Range("I85").Value = 1588.6666
$ws.Range("J85").Value = 3494.6
$ws.Range("K85").Value = 1588.6666
$ws.Range("L85").Value = 3494.6
$ws.Range("M85").Value = -340.6666
$ws.Range("N85").Value = -5990.6
$ws.Range("H126").Value = 5902
$ws.Range("I126").Value = 5902
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 17706
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = -15236
$ws.Range("N126").ClearContents()
$ws.Range("H130").Value = 0
$ws.Range("J130").Value = 0
$ws.Range("L130").Value = 0
$ws.Range("N130").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H48").Value = 0
$ws.Range("I48").Value = 0
$ws.Range("K48").Value = 0
$ws.Range("M48").ClearContents()
$ws.Range("H81").Value = 1396
$ws.Range("J81").Value = 1655
$ws.Range("L81").Value = 3310
$ws.Range("N81").Value = -5432
$ws.Range("H84").Value = 1396
$ws.Range("J84").Value = 1655
$ws.Range("L84").Value = 16550
$ws.Range("N84").Value = -27158
$ws.Range("H122").Value = 87908.78
$ws.Range("I122").Value = 98584.875
$ws.Range("K122").Value = 295754.625
$ws.Range("M122").Value = -293304.625
$ws.Range("H132").Value = 1302.6957
$ws.Range("I132").Value = 1022.0476
$ws.Range("K132").Value = 3066.1428
$ws.Range("M132").Value = -536.1428000000001
